$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.802.82"
$ws.Range("E2").Value = "  +4.35%  "

$ws.Range("D3").Value = "2.421.09"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.68%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +6.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.514"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +2.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +11.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.38"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "18.71"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.00%  "

$ws.Range("E13").Value = "  -1.03%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.91"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.95%  "

$ws.Range("D15").Value = "2.802.83"
$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("D16").Value = "2.440.98"
$ws.Range("E16").Value = "  +3.76%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.832"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +4.80%  "

$ws.Range("D18").Value = "44.639.97"
$ws.Range("E18").Value = "  +4.01%  "

$ws.Range("E19").Value = "  +3.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.36"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.56%  "

$ws.Range("D21").Value = "0.0₃0918"
$ws.Range("E21").Value = "  +3.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.72"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.62"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.29%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.90%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.36"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +3.87%  "

$ws.Range("E28").Value = "  -3.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.52"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.49"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +3.74%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.50"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.36%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +17.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.51"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +11.84%  "

$ws.Range("E34").Value = "  +3.32%  "

$ws.Range("E35").Value = "  +0.22%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0768"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.50"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.19%  "

$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "124.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.36%  "

$ws.Range("E41").Value = "  +2.18%  "

$ws.Range("E42").Value = "  -3.37%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.23%  "

$ws.Range("E44").Value = "  +4.53%  "

$ws.Range("D45").Value = "1.941.78"
$ws.Range("E45").Value = "  +0.79%  "

$ws.Range("E46").Value = "  -1.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.20"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.74"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +15.89%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.54"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.85%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +5.87%  "
